{"js": "// Update the two-digit x two-digit multiplication prompts in the table.\n// Each cell holds a single unique \"NN\u00d7NN=\" text run, so searching the\n// body for each exact old string and replacing it with the new string\n// is safe and unambiguous.\n\nconst replacements = [\n  [\"82\u00d797=\", \"42\u00d748=\"],\n  [\"44\u00d777=\", \"33\u00d786=\"],\n  [\"52\u00d758=\", \"69\u00d731=\"],\n  [\"47\u00d750=\", \"26\u00d777=\"],\n  [\"54\u00d797=\", \"85\u00d774=\"],\n  [\"84\u00d785=\", \"68\u00d781=\"],\n  [\"88\u00d741=\", \"21\u00d737=\"],\n  [\"42\u00d771=\", \"51\u00d793=\"],\n  [\"19\u00d781=\", \"85\u00d726=\"],\n  [\"86\u00d778=\", \"32\u00d740=\"],\n  [\"48\u00d731=\", \"68\u00d724=\"],\n  [\"30\u00d795=\", \"38\u00d718=\"],\n  [\"24\u00d744=\", \"51\u00d744=\"],\n  [\"79\u00d777=\", \"23\u00d734=\"],\n  [\"33\u00d776=\", \"74\u00d736=\"],\n  [\"90\u00d777=\", \"62\u00d793=\"],\n  [\"74\u00d729=\", \"71\u00d787=\"],\n  [\"66\u00d725=\", \"51\u00d797=\"],\n  [\"59\u00d758=\", \"57\u00d792=\"],\n  [\"25\u00d749=\", \"16\u00d797=\"],\n  [\"20\u00d755=\", \"44\u00d752=\"],\n  [\"76\u00d776=\", \"65\u00d767=\"],\n  [\"38\u00d782=\", \"96\u00d780=\"],\n  [\"46\u00d732=\", \"15\u00d779=\"],\n  [\"84\u00d759=\", \"91\u00d770=\"],\n];\n\nconst body = context.document.body;\n\n// Kick off all the searches first, then resolve them in one sync.\nconst pending = replacements.map(([oldText, newText]) => {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  return { results, newText };\n});\n\nawait context.sync();\n\nfor (const { results, newText } of pending) {\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the two-digit x two-digit multiplication prompts in the table.\n# Each cell holds a single unique \"NN\u00d7NN=\" run, so a plain Find/Replace\n# (MatchCase + MatchWholeWord-ish exact text) for each pair is safe and\n# unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"82\u00d797=\", \"42\u00d748=\"),\n    @(\"44\u00d777=\", \"33\u00d786=\"),\n    @(\"52\u00d758=\", \"69\u00d731=\"),\n    @(\"47\u00d750=\", \"26\u00d777=\"),\n    @(\"54\u00d797=\", \"85\u00d774=\"),\n    @(\"84\u00d785=\", \"68\u00d781=\"),\n    @(\"88\u00d741=\", \"21\u00d737=\"),\n    @(\"42\u00d771=\", \"51\u00d793=\"),\n    @(\"19\u00d781=\", \"85\u00d726=\"),\n    @(\"86\u00d778=\", \"32\u00d740=\"),\n    @(\"48\u00d731=\", \"68\u00d724=\"),\n    @(\"30\u00d795=\", \"38\u00d718=\"),\n    @(\"24\u00d744=\", \"51\u00d744=\"),\n    @(\"79\u00d777=\", \"23\u00d734=\"),\n    @(\"33\u00d776=\", \"74\u00d736=\"),\n    @(\"90\u00d777=\", \"62\u00d793=\"),\n    @(\"74\u00d729=\", \"71\u00d787=\"),\n    @(\"66\u00d725=\", \"51\u00d797=\"),\n    @(\"59\u00d758=\", \"57\u00d792=\"),\n    @(\"25\u00d749=\", \"16\u00d797=\"),\n    @(\"20\u00d755=\", \"44\u00d752=\"),\n    @(\"76\u00d776=\", \"65\u00d767=\"),\n    @(\"38\u00d782=\", \"96\u00d780=\"),\n    @(\"46\u00d732=\", \"15\u00d779=\"),\n    @(\"84\u00d759=\", \"91\u00d770=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
